$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = 5.100300000000002
$ws.Range("A9").Value = -21.87670000000002
$ws.Range("B12").Value = 5.331099999999996
$ws.Range("D13").Value = -8.492500000000001
$ws.Range("C15").Value = -13.80189999999999
$ws.Range("D16").Value = -8.643700000000006
$ws.Range("A18").Value = -22.13660000000001
$ws.Range("A20").Value = -19.91909999999998
$ws.Range("D20").Value = -7.144699999999998
$ws.Range("D24").Value = -7.271900000000003
$ws.Range("B26").Value = 4.275900000000006
$ws.Range("A27").Value = -21.64369999999998
$ws.Range("B27").Value = 5.411799999999999
$ws.Range("B29").Value = 4.825299999999997
$ws.Range("B37").Value = 8.660300000000012
$ws.Range("B38").Value = 4.897600000000003
$ws.Range("C38").Value = -12.1484
$ws.Range("D39").Value = -7.491600000000003
$ws.Range("C44").Value = -13.34389999999999
$ws.Range("D48").Value = -7.234899999999998
$ws.Range("B51").Value = 6.003700000000002
$ws.Range("C51").Value = -11.9312
$ws.Range("D52").Value = -7.763500000000001
$ws.Range("B55").Value = 5.038599999999997
$ws.Range("D56").Value = -7.846399999999996
$ws.Range("C57").Value = -14.07539999999999
$ws.Range("C63").Value = -11.22090000000001
$ws.Range("A69").Value = -21.6523
$ws.Range("B69").Value = 5.377399999999998
$ws.Range("B70").Value = 6.388000000000006
$ws.Range("C70").Value = -11.7182
$ws.Range("A76").Value = -19.47569999999999
$ws.Range("A82").Value = -21.90010000000001
$ws.Range("B83").Value = 6.110300000000002
$ws.Range("D84").Value = -8.951099999999999
$ws.Range("C99").Value = -12.8804
$ws.Range("D100").Value = -8.279500000000004
$ws.Range("D101").Value = -7.801099999999996
$ws.Range("B102").Value = 7.945800000000007
